$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds serial date 45190 (2023-09-21) for every
# data row (rows 2-151). Update it to 45192 (2023-09-23).
for ($r = 2; $r -le 151; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45190) {
        $cell.Value2 = 45192
    }
}
